# Applies the automated daily data refresh to the stock data sheet.
# Updates intraday / snapshot derived metrics (Open, Gap %, VWAP %, and the
# High/Low pairs for the 1m and 5m windows) for each ticker row (2-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - DWTX
$ws.Range("H2").Value = 8.41
$ws.Range("O2").Value = 52.36
$ws.Range("Q2").Value = -9.01
$ws.Range("AJ2").Value = 8.550000000000001
$ws.Range("AK2").Value = 8.19
$ws.Range("AN2").Value = 8.550000000000001
$ws.Range("AO2").Value = 8.19

# Row 3 - JFB
$ws.Range("H3").Value = 11.21
$ws.Range("O3").Value = 63.65
$ws.Range("P3").Value = 11.56
$ws.Range("Q3").Value = 3.08
$ws.Range("AJ3").Value = 11.21
$ws.Range("AK3").Value = 10.65
$ws.Range("AN3").Value = 11.21
$ws.Range("AO3").Value = 10.65

# Row 4 - MSS
$ws.Range("H4").Value = 3.8
$ws.Range("O4").Value = 297.91
$ws.Range("Q4").Value = -56.89
$ws.Range("T4").Value = 3.96
$ws.Range("X4").Value = 3.96
$ws.Range("AB4").Value = 3.96
$ws.Range("AF4").Value = 3.96
$ws.Range("AJ4").Value = 3.96
$ws.Range("AN4").Value = 3.96

# Row 5 - POAI
$ws.Range("H5").Value = 2.13
$ws.Range("O5").Value = 175.48
$ws.Range("Q5").Value = -30.27
$ws.Range("T5").Value = 2.13
$ws.Range("X5").Value = 2.13
$ws.Range("AB5").Value = 2.13
$ws.Range("AF5").Value = 2.13
$ws.Range("AJ5").Value = 2.13
$ws.Range("AK5").Value = 1.9
$ws.Range("AN5").Value = 2.13
$ws.Range("AO5").Value = 1.9

# Row 6 - ZURA
$ws.Range("H6").Value = 3.59
$ws.Range("O6").Value = 48.76
$ws.Range("Q6").Value = -9.550000000000001
